# daily auto push: 2026-01-22 02:33 UTC
#
# Insert a new daily-log row for 2026/01/22 (Thursday) ahead of the
# 2026/12/29 block, shifting every following row down by one, exactly
# like Excel's normal "insert row" behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 675 (and everything below it) down by one row.
$ws.Rows(675).Insert()

# Column A holds a date-look-alike string ("2026/01/22"). A plain
# assignment would make Excel auto-convert it to a date serial number,
# so a leading apostrophe is used to force text entry, and the style is
# then reset back to Normal so the quote-prefix formatting doesn't stick
# to the cell (matching the unstyled look of every other data row).
$ws.Range("A675").Value = "'2026/01/22"
$ws.Range("A675").Style = "Normal"

$ws.Range("B675").Value = "木"
$ws.Range("C675").Value = 10
$ws.Range("D675").Value = 17
